# Auto-generated Excel COM-interop script
# Applies numeric value updates (scheduled price-data refresh) across all 8 sheets
# as described by the OOXML diff for Sheets/Hyperion_Profits.xlsx.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 261.25  # H9: was 310.1111
$ws.Cells.Item(9, 9).Value = 284.2857  # I9: was 336.375
$ws.Cells.Item(9, 11).Value = 284.2857  # K9: was 336.375
$ws.Cells.Item(9, 13).Value = -115.2857  # M9: was -167.375
$ws.Cells.Item(32, 8).Value = 4292.9614  # H32: was 4386.64
$ws.Cells.Item(32, 9).Value = 1984  # I32: was 1994.5
$ws.Cells.Item(32, 10).Value = 4485.375  # J32: was 4594.6523
$ws.Cells.Item(32, 11).Value = 1984  # K32: was 1994.5
$ws.Cells.Item(32, 12).Value = 4485.375  # L32: was 4594.6523
$ws.Cells.Item(32, 13).Value = -1658  # M32: was -1668.5
$ws.Cells.Item(32, 14).Value = -5137.375  # N32: was -5246.6523
$ws.Cells.Item(53, 8).Value = 6877.933  # H53: was 6448.375
$ws.Cells.Item(53, 9).Value = 274.2857  # I53: was 240.125
$ws.Cells.Item(53, 10).Value = 12656.125  # J53: was 12656.625
$ws.Cells.Item(53, 11).Value = 274.2857  # K53: was 240.125
$ws.Cells.Item(53, 12).Value = 12656.125  # L53: was 12656.625
$ws.Cells.Item(53, 13).Value = 362.7143  # M53: was 396.875
$ws.Cells.Item(53, 14).Value = -13930.125  # N53: was -13930.625
$ws.Cells.Item(62, 8).Value = 5371.7827  # H62: was 5376.087
$ws.Cells.Item(62, 9).Value = 3567.7368  # I62: was 3666
$ws.Cells.Item(62, 10).Value = 13941  # J62: was 11532.4
$ws.Cells.Item(62, 11).Value = 3567.7368  # K62: was 3666
$ws.Cells.Item(62, 12).Value = 13941  # L62: was 11532.4
$ws.Cells.Item(62, 13).Value = -2943.7368  # M62: was -3042
$ws.Cells.Item(62, 14).Value = -15189  # N62: was -12780.4
$ws.Cells.Item(65, 8).Value = 5371.7827  # H65: was 5376.087
$ws.Cells.Item(65, 9).Value = 3567.7368  # I65: was 3666
$ws.Cells.Item(65, 10).Value = 13941  # J65: was 11532.4
$ws.Cells.Item(65, 11).Value = 17838.684  # K65: was 18330
$ws.Cells.Item(65, 12).Value = 69705  # L65: was 57662
$ws.Cells.Item(65, 13).Value = -14718.684  # M65: was -15210
$ws.Cells.Item(65, 14).Value = -75945  # N65: was -63902
$ws.Cells.Item(86, 8).Value = 1811.9286  # H86: was 1741.6875
$ws.Cells.Item(86, 9).Value = 1967  # I86: was 1872.5
$ws.Cells.Item(86, 10).Value = 1725.7778  # J86: was 1663.2
$ws.Cells.Item(86, 11).Value = 1967  # K86: was 1872.5
$ws.Cells.Item(86, 12).Value = 1725.7778  # L86: was 1663.2
$ws.Cells.Item(86, 13).Value = -844  # M86: was -749.5
$ws.Cells.Item(86, 14).Value = -3971.7778  # N86: was -3909.2
$ws.Cells.Item(89, 8).Value = 1811.9286  # H89: was 1741.6875
$ws.Cells.Item(89, 9).Value = 1967  # I89: was 1872.5
$ws.Cells.Item(89, 10).Value = 1725.7778  # J89: was 1663.2
$ws.Cells.Item(89, 11).Value = 9835  # K89: was 9362.5
$ws.Cells.Item(89, 12).Value = 8628.889000000001  # L89: was 8316
$ws.Cells.Item(89, 13).Value = -4219  # M89: was -3746.5
$ws.Cells.Item(89, 14).Value = -19860.889  # N89: was -19548
$ws.Cells.Item(99, 8).Value = 311.9  # H99: was 369.27274
$ws.Cells.Item(99, 9).Value = 248.66667  # I99: was 256.8889
$ws.Cells.Item(99, 10).Value = 881  # J99: was 875
$ws.Cells.Item(99, 11).Value = 746.00001  # K99: was 770.6667
$ws.Cells.Item(99, 12).Value = 2643  # L99: was 2625
$ws.Cells.Item(99, 13).Value = 751.99999  # M99: was 727.3333
$ws.Cells.Item(99, 14).Value = -5639  # N99: was -5621
$ws.Cells.Item(101, 8).Value = 33334810  # H101: was 38463100
$ws.Cells.Item(101, 9).Value = 62500260  # I101: was 71428750
$ws.Cells.Item(101, 10).Value = 2864.1428  # J101: was 3166.3333
$ws.Cells.Item(101, 11).Value = 187500780  # K101: was 214286250
$ws.Cells.Item(101, 12).Value = 8592.428400000001  # L101: was 9498.999899999999
$ws.Cells.Item(101, 13).Value = -187499158  # M101: was -214284628
$ws.Cells.Item(101, 14).Value = -11836.4284  # N101: was -12742.9999
$ws.Cells.Item(113, 8).Value = 6645.609  # H113: was 6485.522
$ws.Cells.Item(113, 9).Value = 6585.5557  # I113: was 6593.778
$ws.Cells.Item(113, 10).Value = 6684.2144  # J113: was 6415.9287
$ws.Cells.Item(113, 11).Value = 6585.5557  # K113: was 6593.778
$ws.Cells.Item(113, 12).Value = 6684.2144  # L113: was 6415.9287
$ws.Cells.Item(113, 13).Value = -3331.5557  # M113: was -3339.778
$ws.Cells.Item(113, 14).Value = -13192.2144  # N113: was -12923.9287
$ws.Cells.Item(127, 8).Value = 1072.3334  # H127: was 776.6
$ws.Cells.Item(127, 9).Value = 1000  # I127: was 695.75
$ws.Cells.Item(127, 10).Value = 1217  # J127: was 1100
$ws.Cells.Item(127, 11).Value = 3000  # K127: was 2087.25
$ws.Cells.Item(127, 12).Value = 3651  # L127: was 3300
$ws.Cells.Item(127, 13).Value = 1960  # M127: was 2872.75
$ws.Cells.Item(127, 14).Value = -13571  # N127: was -13220
$ws.Cells.Item(129, 8).Value = 23810580  # H129: was 25001068
$ws.Cells.Item(129, 9).Value = 26316748  # I129: was 26316756
$ws.Cells.Item(129, 10).Value = 2000  # J129: was 3000
$ws.Cells.Item(129, 11).Value = 78950244  # K129: was 78950268
$ws.Cells.Item(129, 12).Value = 6000  # L129: was 9000
$ws.Cells.Item(129, 13).Value = -78945244  # M129: was -78945268
$ws.Cells.Item(129, 14).Value = -16000  # N129: was -19000

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 280  # H5: was 189.5
$ws.Cells.Item(5, 9).Value = 287.5  # I5: was 189.5
$ws.Cells.Item(5, 10).Value = 250  # J5: was 0
$ws.Cells.Item(5, 11).Value = 287.5  # K5: was 189.5
$ws.Cells.Item(5, 12).Value = 250  # L5: was 0
$ws.Cells.Item(5, 13).Value = -175.5  # M5: was -77.5
$ws.Cells.Item(5, 14).Value = -474  # N5: was None
$ws.Cells.Item(32, 8).Value = 2521.92  # H32: was 2521.93
$ws.Cells.Item(32, 9).Value = 1663.875  # I32: was 1663.8864
$ws.Cells.Item(32, 11).Value = 1663.875  # K32: was 1663.8864
$ws.Cells.Item(32, 13).Value = -1376.875  # M32: was -1376.8864
$ws.Cells.Item(49, 8).Value = 10000  # H49: was 9520
$ws.Cells.Item(49, 10).Value = 10000  # J49: was 9520
$ws.Cells.Item(49, 12).Value = 10000  # L49: was 9520
$ws.Cells.Item(49, 14).Value = -10520  # N49: was -10040
$ws.Cells.Item(61, 8).Value = 2269.1  # H61: was 2986.9092
$ws.Cells.Item(61, 9).Value = 1835.5  # I61: was 2605.25
$ws.Cells.Item(61, 10).Value = 4003.5  # J61: was 4004.6667
$ws.Cells.Item(61, 11).Value = 1835.5  # K61: was 2605.25
$ws.Cells.Item(61, 12).Value = 4003.5  # L61: was 4004.6667
$ws.Cells.Item(61, 13).Value = -1623.5  # M61: was -2393.25
$ws.Cells.Item(61, 14).Value = -4427.5  # N61: was -4428.6667
$ws.Cells.Item(132, 8).Value = 2694.0908  # H132: was 2867.6843
$ws.Cells.Item(132, 9).Value = 3086.889  # I132: was 3833
$ws.Cells.Item(132, 11).Value = 9260.667000000001  # K132: was 11499
$ws.Cells.Item(132, 13).Value = -6730.667000000001  # M132: was -8969
$ws.Cells.Item(136, 8).Value = 2269.1  # H136: was 2986.9092
$ws.Cells.Item(136, 9).Value = 1835.5  # I136: was 2605.25
$ws.Cells.Item(136, 10).Value = 4003.5  # J136: was 4004.6667
$ws.Cells.Item(136, 11).Value = 5506.5  # K136: was 7815.75
$ws.Cells.Item(136, 12).Value = 12010.5  # L136: was 12014.0001
$ws.Cells.Item(136, 13).Value = -2956.5  # M136: was -5265.75
$ws.Cells.Item(136, 14).Value = -17110.5  # N136: was -17114.0001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 280  # H4: was 189.5
$ws.Cells.Item(4, 9).Value = 287.5  # I4: was 189.5
$ws.Cells.Item(4, 10).Value = 250  # J4: was 0
$ws.Cells.Item(4, 11).Value = 287.5  # K4: was 189.5
$ws.Cells.Item(4, 12).Value = 250  # L4: was 0
$ws.Cells.Item(4, 13).Value = -172.5  # M4: was -74.5
$ws.Cells.Item(4, 14).Value = -480  # N4: was None
$ws.Cells.Item(86, 8).Value = 3453462.8  # H86: was 4172822
$ws.Cells.Item(86, 9).Value = 6257795.5  # I86: was 8343449
$ws.Cells.Item(86, 10).Value = 1976.4615  # J86: was 2195.3333
$ws.Cells.Item(86, 11).Value = 6257795.5  # K86: was 8343449
$ws.Cells.Item(86, 12).Value = 1976.4615  # L86: was 2195.3333
$ws.Cells.Item(86, 13).Value = -6256672.5  # M86: was -8342326
$ws.Cells.Item(86, 14).Value = -4222.461499999999  # N86: was -4441.3333
$ws.Cells.Item(89, 8).Value = 3453462.8  # H89: was 4172822
$ws.Cells.Item(89, 9).Value = 6257795.5  # I89: was 8343449
$ws.Cells.Item(89, 10).Value = 1976.4615  # J89: was 2195.3333
$ws.Cells.Item(89, 11).Value = 31288977.5  # K89: was 41717245
$ws.Cells.Item(89, 12).Value = 9882.307499999999  # L89: was 10976.6665
$ws.Cells.Item(89, 13).Value = -31283361.5  # M89: was -41711629
$ws.Cells.Item(89, 14).Value = -21114.3075  # N89: was -22208.6665
$ws.Cells.Item(104, 8).Value = 54984  # H104: was 0
$ws.Cells.Item(104, 10).Value = 54984  # J104: was 0
$ws.Cells.Item(104, 12).Value = 54984  # L104: was 0
$ws.Cells.Item(104, 14).Value = -61972  # N104: was None
$ws.Cells.Item(109, 8).Value = 44342  # H109: was 68682
$ws.Cells.Item(109, 10).Value = 44342  # J109: was 68682
$ws.Cells.Item(109, 12).Value = 44342  # L109: was 68682
$ws.Cells.Item(109, 14).Value = -47116  # N109: was -71456
$ws.Cells.Item(139, 8).Value = 114294.6  # H139: was 116534.6
$ws.Cells.Item(139, 10).Value = 147249  # J139: was 152849
$ws.Cells.Item(139, 12).Value = 147249  # L139: was 152849
$ws.Cells.Item(139, 14).Value = -157529  # N139: was -163129

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 371.57895  # H7: was 403.35294
$ws.Cells.Item(7, 9).Value = 243.3077  # I7: was 263.33334
$ws.Cells.Item(7, 10).Value = 649.5  # J7: was 739.4
$ws.Cells.Item(7, 11).Value = 243.3077  # K7: was 263.33334
$ws.Cells.Item(7, 12).Value = 649.5  # L7: was 739.4
$ws.Cells.Item(7, 13).Value = -130.3077  # M7: was -150.33334
$ws.Cells.Item(7, 14).Value = -875.5  # N7: was -965.4
$ws.Cells.Item(31, 8).Value = 5124.5  # H31: was 5156.4287
$ws.Cells.Item(31, 10).Value = 5700.5  # J31: was 6500
$ws.Cells.Item(31, 12).Value = 5700.5  # L31: was 6500
$ws.Cells.Item(31, 14).Value = -6290.5  # N31: was -7090
$ws.Cells.Item(34, 8).Value = 5124.5  # H34: was 5156.4287
$ws.Cells.Item(34, 10).Value = 5700.5  # J34: was 6500
$ws.Cells.Item(34, 12).Value = 5700.5  # L34: was 6500
$ws.Cells.Item(34, 14).Value = -6104.5  # N34: was -6904
$ws.Cells.Item(99, 8).Value = 4514.625  # H99: was 4642.2666
$ws.Cells.Item(99, 9).Value = 4078.9  # I99: was 4243.222
$ws.Cells.Item(99, 11).Value = 4078.9  # K99: was 4243.222
$ws.Cells.Item(99, 13).Value = -2580.9  # M99: was -2745.222
$ws.Cells.Item(126, 8).Value = 4514.625  # H126: was 4642.2666
$ws.Cells.Item(126, 9).Value = 4078.9  # I126: was 4243.222
$ws.Cells.Item(126, 11).Value = 12236.7  # K126: was 12729.666
$ws.Cells.Item(126, 13).Value = -9766.700000000001  # M126: was -10259.666
$ws.Cells.Item(132, 8).Value = 2985.3572  # H132: was 2443.641
$ws.Cells.Item(132, 9).Value = 2837.05  # I132: was 2277.6897
$ws.Cells.Item(132, 10).Value = 3356.125  # J132: was 2924.9
$ws.Cells.Item(132, 11).Value = 8511.150000000001  # K132: was 6833.0691
$ws.Cells.Item(132, 12).Value = 10068.375  # L132: was 8774.700000000001
$ws.Cells.Item(132, 13).Value = -5981.150000000001  # M132: was -4303.0691
$ws.Cells.Item(132, 14).Value = -15128.375  # N132: was -13834.7
$ws.Cells.Item(133, 8).Value = 30326  # H133: was 37500
$ws.Cells.Item(133, 10).Value = 30326  # J133: was 37500
$ws.Cells.Item(133, 12).Value = 30326  # L133: was 37500
$ws.Cells.Item(133, 14).Value = -35386  # N133: was -42560

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 208.48936  # H2: was 224.06818
$ws.Cells.Item(2, 9).Value = 100.92593  # I2: was 112.291664
$ws.Cells.Item(2, 10).Value = 353.7  # J2: was 358.2
$ws.Cells.Item(2, 11).Value = 605.55558  # K2: was 673.749984
$ws.Cells.Item(2, 12).Value = 2122.2  # L2: was 2149.2
$ws.Cells.Item(2, 13).Value = -492.55558  # M2: was -560.749984
$ws.Cells.Item(2, 14).Value = -2348.2  # N2: was -2375.2
$ws.Cells.Item(70, 8).Value = 670  # H70: was 760
$ws.Cells.Item(70, 9).Value = 662.5  # I70: was 700
$ws.Cells.Item(70, 10).Value = 700  # J70: was 850
$ws.Cells.Item(70, 11).Value = 1987.5  # K70: was 2100
$ws.Cells.Item(70, 12).Value = 2100  # L70: was 2550
$ws.Cells.Item(70, 13).Value = -1672.5  # M70: was -1785
$ws.Cells.Item(70, 14).Value = -2730  # N70: was -3180
$ws.Cells.Item(73, 8).Value = 670  # H73: was 760
$ws.Cells.Item(73, 9).Value = 662.5  # I73: was 700
$ws.Cells.Item(73, 10).Value = 700  # J73: was 850
$ws.Cells.Item(73, 11).Value = 1987.5  # K73: was 2100
$ws.Cells.Item(73, 12).Value = 2100  # L73: was 2550
$ws.Cells.Item(73, 13).Value = -895.5  # M73: was -1008
$ws.Cells.Item(73, 14).Value = -4284  # N73: was -4734
$ws.Cells.Item(128, 8).Value = 199193.6  # H128: was 199159.5
$ws.Cells.Item(128, 9).Value = 199193.6  # I128: was 199159.5
$ws.Cells.Item(128, 11).Value = 597580.8  # K128: was 597478.5
$ws.Cells.Item(128, 13).Value = -592600.8  # M128: was -592498.5
$ws.Cells.Item(136, 8).Value = 2100  # H136: was 2233.1667
$ws.Cells.Item(136, 9).Value = 2100  # I136: was 2233.1667
$ws.Cells.Item(136, 11).Value = 6300  # K136: was 6699.500100000001
$ws.Cells.Item(136, 13).Value = -1200  # M136: was -1599.500100000001
$ws.Cells.Item(138, 8).Value = 14631.0625  # H138: was 15549.866
$ws.Cells.Item(138, 9).Value = 14631.0625  # I138: was 15549.866
$ws.Cells.Item(138, 11).Value = 43893.1875  # K138: was 46649.598
$ws.Cells.Item(138, 13).Value = -38753.1875  # M138: was -41509.598

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 68303  # H122: was 70797.81
$ws.Cells.Item(122, 9).Value = 86072  # I122: was 82371.35000000001
$ws.Cells.Item(122, 10).Value = 3150  # J122: was 4250
$ws.Cells.Item(122, 11).Value = 258216  # K122: was 247114.05
$ws.Cells.Item(122, 12).Value = 9450  # L122: was 12750
$ws.Cells.Item(122, 13).Value = -255766  # M122: was -244664.05
$ws.Cells.Item(122, 14).Value = -14350  # N122: was -17650
$ws.Cells.Item(137, 8).Value = 122498  # H137: was 127000
$ws.Cells.Item(137, 10).Value = 122498  # J137: was 127000
$ws.Cells.Item(137, 12).Value = 122498  # L137: was 127000
$ws.Cells.Item(137, 14).Value = -132698  # N137: was -137200

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 15875386  # H93: was 15153821
$ws.Cells.Item(93, 9).Value = 19610206  # I93: was 18520804
$ws.Cells.Item(93, 11).Value = 19610206  # K93: was 18520804
$ws.Cells.Item(93, 13).Value = -19608958  # M93: was -18519556
$ws.Cells.Item(100, 8).Value = 38162.465  # H100: was 39243.93
$ws.Cells.Item(100, 9).Value = 4314.96  # I100: was 4322.96
$ws.Cells.Item(100, 10).Value = 207400  # J100: was 257500
$ws.Cells.Item(100, 11).Value = 4314.96  # K100: was 4322.96
$ws.Cells.Item(100, 12).Value = 207400  # L100: was 257500
$ws.Cells.Item(100, 13).Value = -3773.96  # M100: was -3781.96
$ws.Cells.Item(100, 14).Value = -208482  # N100: was -258582
$ws.Cells.Item(109, 8).Value = 49981.5  # H109: was 49982
$ws.Cells.Item(109, 10).Value = 49981.5  # J109: was 49982
$ws.Cells.Item(109, 12).Value = 49981.5  # L109: was 49982
$ws.Cells.Item(109, 14).Value = -52755.5  # N109: was -52756
$ws.Cells.Item(132, 8).Value = 5668.511  # H132: was 6290.775
$ws.Cells.Item(132, 9).Value = 6012.6284  # I132: was 6930.6895
$ws.Cells.Item(132, 10).Value = 4464.1  # J132: was 4603.727
$ws.Cells.Item(132, 11).Value = 18037.8852  # K132: was 20792.0685
$ws.Cells.Item(132, 12).Value = 13392.3  # L132: was 13811.181
$ws.Cells.Item(132, 13).Value = -15507.8852  # M132: was -18262.0685
$ws.Cells.Item(132, 14).Value = -18452.3  # N132: was -18871.181
$ws.Cells.Item(136, 8).Value = 77992.71000000001  # H136: was 77995.78999999999
$ws.Cells.Item(136, 9).Value = 203979.6  # I136: was 254152
$ws.Cells.Item(136, 10).Value = 8000  # J136: was 7533.3
$ws.Cells.Item(136, 11).Value = 611938.8  # K136: was 762456
$ws.Cells.Item(136, 12).Value = 24000  # L136: was 22599.9
$ws.Cells.Item(136, 13).Value = -609388.8  # M136: was -759906
$ws.Cells.Item(136, 14).Value = -29100  # N136: was -27699.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(80, 8).Value = 53859.6  # H80: was 59824.75
$ws.Cells.Item(80, 9).Value = 0  # I80: was 30000
$ws.Cells.Item(80, 10).Value = 53859.6  # J80: was 69766.336
$ws.Cells.Item(80, 11).Value = 0  # K80: was 30000
$ws.Cells.Item(80, 12).Value = 53859.6  # L80: was 69766.336
$ws.Cells.Item(80, 13).ClearContents()  # M80: was -29002
$ws.Cells.Item(80, 14).Value = -55855.6  # N80: was -71762.336
$ws.Cells.Item(83, 8).Value = 53859.6  # H83: was 59824.75
$ws.Cells.Item(83, 9).Value = 0  # I83: was 30000
$ws.Cells.Item(83, 10).Value = 53859.6  # J83: was 69766.336
$ws.Cells.Item(83, 11).Value = 0  # K83: was 90000
$ws.Cells.Item(83, 12).Value = 161578.8  # L83: was 209299.008
$ws.Cells.Item(83, 13).ClearContents()  # M83: was -85008
$ws.Cells.Item(83, 14).Value = -171562.8  # N83: was -219283.008
$ws.Cells.Item(109, 8).Value = 47992.5  # H109: was 47994
$ws.Cells.Item(109, 10).Value = 47992.5  # J109: was 47994
$ws.Cells.Item(109, 12).Value = 47992.5  # L109: was 47994
$ws.Cells.Item(109, 14).Value = -50766.5  # N109: was -50768
$ws.Cells.Item(132, 8).Value = 28575646  # H132: was 30334442
$ws.Cells.Item(132, 9).Value = 33337892  # I132: was 38466716
$ws.Cells.Item(132, 10).Value = 2169.8  # J132: was 128851.29
$ws.Cells.Item(132, 11).Value = 100013676  # K132: was 115400148
$ws.Cells.Item(132, 12).Value = 6509.400000000001  # L132: was 386553.87
$ws.Cells.Item(132, 13).Value = -100011146  # M132: was -115397618
$ws.Cells.Item(132, 14).Value = -11569.4  # N132: was -391613.87
